$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the oldest data block (rows 1164-1166) down to brand new rows
#        1167-1169, since a new weekly entry pushes it out of the rolling window.
$src = $ws.Range("A1164:T1166")
$dst = $ws.Range("A1167:T1169")
$src.Copy($dst)

# --- 2. Read the existing 15 data blocks (rows 1122-1166, 3 rows each) so we can
#        shift them down by one block (each block takes on the values that used
#        to belong to the block above it).
$blocks = @()
for ($b = 0; $b -lt 15; $b++) {
    $baseRow = 1122 + ($b * 3)
    $rowsData = @()
    for ($j = 0; $j -lt 3; $j++) {
        $r = $baseRow + $j
        $d = $ws.Cells.Item($r, 4).Value2
        $m = $ws.Cells.Item($r, 13).Value2
        $n = $ws.Cells.Item($r, 14).Value2
        $o = $ws.Cells.Item($r, 15).Value2
        $p = $ws.Cells.Item($r, 16).Value2
        $s = $ws.Cells.Item($r, 19).Value2
        $rowsData += ,@($d, $m, $n, $o, $p, $s)
    }
    $blocks += ,$rowsData
}

# --- 3. Brand new weekly data block (the newest week) that now occupies rows
#        1122-1124, pushing every other block down by one position.
$newBlock = @(
    @(45075, 80,  14000, 14000, 14000, 700),
    @(45075, 120, 16000, 16000, 16000, 800),
    @(45075, 120, 17000, 17000, 17000, 850)
)

# --- 4. Build the final sequence of 15 blocks for rows 1122-1166: the new block
#        followed by the first 14 of the original blocks (the 15th/oldest one was
#        already copied out to rows 1167-1169 above).
$finalBlocks = @()
$finalBlocks += ,$newBlock
for ($b = 0; $b -lt 14; $b++) {
    $finalBlocks += ,$blocks[$b]
}

# --- 5. Write the final blocks back into rows 1122-1166.
for ($b = 0; $b -lt 15; $b++) {
    $baseRow = 1122 + ($b * 3)
    $rowsData = $finalBlocks[$b]
    for ($j = 0; $j -lt 3; $j++) {
        $r = $baseRow + $j
        $vals = $rowsData[$j]
        $ws.Cells.Item($r, 4).Value = $vals[0]
        $ws.Cells.Item($r, 13).Value = $vals[1]
        $ws.Cells.Item($r, 14).Value = $vals[2]
        $ws.Cells.Item($r, 15).Value = $vals[3]
        $ws.Cells.Item($r, 16).Value = $vals[4]
        $ws.Cells.Item($r, 19).Value = $vals[5]
    }
}
